$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999465098655,
    0.99827606596830698,
    0.99260386160420189,
    0.99425943816268803,
    0.98271339805427194,
    0.95857168674844773,
    0.9519742358858575,
    0.94284274751499852,
    0.93109531747852214,
    0.92026554619628298,
    0.91866923426905767,
    0.91588910126035072,
    0.90460096102039123,
    0.90043335614479725,
    0.89784173707118242,
    0.89533515043609158,
    0.89162711863905919,
    0.89051819418988321,
    0.99722284377789627,
    0.99010581337011105,
    0.98870732036345521,
    0.98744281479215545,
    0.98556961252467623,
    0.97254929785723121,
    0.96609239354003695,
    0.95645981192247387,
    0.9537650351627992,
    0.94443511707994809,
    0.93817631325745188,
    0.93530636244219045,
    0.93113530735709382,
    0.92945601361753383,
    0.92893600128157727
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
